$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.508.43"
$ws.Range("E2").Value = "  -0.86%  "
$ws.Range("D3").Value = "2.161.89"
$ws.Range("E3").Value = "  -2.63%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.63"
$ws.Range("E5").Value = "  -1.66%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.605"
$ws.Range("E6").Value = "  -3.24%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "71.41"
$ws.Range("E7").Value = "  -1.94%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  -3.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.56"
$ws.Range("E10").Value = "  -6.35%  "
$ws.Range("E11").Value = "  -4.49%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "53.82"
$ws.Range("E12").Value = "  -4.93%  "
$ws.Range("E13").Value = "  -3.66%  "
$ws.Range("E14").Value = "  -4.16%  "
$ws.Range("D15").Value = "2.484.47"
$ws.Range("E15").Value = "  -2.61%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.14"
$ws.Range("E16").Value = "  -0.97%  "
$ws.Range("D17").Value = "2.154.25"
$ws.Range("E17").Value = "  -2.96%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.778"
$ws.Range("E18").Value = "  -6.90%  "
$ws.Range("D19").Value = "41.404.43"
$ws.Range("E19").Value = "  -0.69%  "
$ws.Range("E20").Value = "  -3.98%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "69.57"
$ws.Range("E21").Value = "  -4.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.76"
$ws.Range("E22").Value = "  -6.70%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.95"
$ws.Range("E23").Value = "  -9.71%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "227.55"
$ws.Range("E24").Value = "  -0.86%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.99"
$ws.Range("E25").Value = "  -4.25%  "
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.65"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.27"
$ws.Range("E28").Value = "  -9.65%  "
$ws.Range("E29").Value = "  -4.41%  "
$ws.Range("E30").Value = "  -0.92%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "171.51"
$ws.Range("E31").Value = "  +2.74%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.72"
$ws.Range("E32").Value = "  -3.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "32.92"
$ws.Range("E33").Value = "  +9.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0769"
$ws.Range("E34").Value = "  -3.40%  "
$ws.Range("E35").Value = "  -7.79%  "
$ws.Range("E36").Value = "  -4.02%  "
$ws.Range("E37").Value = "  -3.26%  "
$ws.Range("E38").Value = "  -1.67%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0299"
$ws.Range("E39").Value = "  -1.33%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.92"
$ws.Range("E40").Value = "  -10.77%  "
$ws.Range("E41").Value = "  -2.85%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.32"
$ws.Range("E42").Value = "  -5.44%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "58.62"
$ws.Range("E43").Value = "  -9.10%  "
$ws.Range("E44").Value = "  -4.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.37"
$ws.Range("E45").Value = "  -3.97%  "
$ws.Range("E46").Value = "  -4.95%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "96.05"
$ws.Range("E47").Value = "  -6.97%  "
$ws.Range("E48").Value = "  -3.46%  "
$ws.Range("E49").Value = "  -4.93%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.16"
$ws.Range("E50").Value = "  -7.50%  "
$ws.Range("E51").Value = "  -2.53%  "
